$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 360
$ws.Range("F2").Value = 45992
$ws.Range("G2").Value = 30865
$ws.Range("H2").Value = 46055
$ws.Range("E3").Value = 30803
$ws.Range("F3").Value = 45992
$ws.Range("E4").Value = 30803
$ws.Range("F4").Value = 45992
$ws.Range("G4").Value = 30865
$ws.Range("H4").Value = 46055
$ws.Range("E5").Value = 30803
$ws.Range("F5").Value = 45992
$ws.Range("G5").Value = 30865
$ws.Range("H5").Value = 46055
$ws.Range("G6").Value = 30865
$ws.Range("H6").Value = 46055
$ws.Range("E7").Value = 30773
$ws.Range("F7").Value = 45962
$ws.Range("G7").Value = 30865
$ws.Range("H7").Value = 46055
$ws.Range("D8").Value = 436
$ws.Range("H8").Value = 46055
$ws.Range("G9").Value = 30865
$ws.Range("H9").Value = 46055
$ws.Range("E10").Value = 30803
$ws.Range("F10").Value = 45992
$ws.Range("G10").Value = 30865
$ws.Range("H10").Value = 46055
$ws.Range("E11").Value = 30773
$ws.Range("F11").Value = 45962
$ws.Range("G11").Value = 30865
$ws.Range("H11").Value = 46055
$ws.Range("C12").Value = 397
$ws.Range("D12").Value = 378
$ws.Range("F12").Value = 45992
$ws.Range("H12").Value = 46055
$ws.Range("C13").Value = 493
$ws.Range("F13").Value = 45992
$ws.Range("G13").Value = 30865
$ws.Range("H13").Value = 46055
$ws.Range("C14").Value = 448
$ws.Range("D14").Value = 422
$ws.Range("F14").Value = 45962
$ws.Range("H14").Value = 46055
$ws.Range("C15").Value = 409
$ws.Range("F15").Value = 45962
$ws.Range("G15").Value = 30834
$ws.Range("H15").Value = 46055
$ws.Range("D16").Value = 436
$ws.Range("H16").Value = 46055
$ws.Range("C17").Value = 396
$ws.Range("D17").Value = 420
$ws.Range("F17").Value = 45992
$ws.Range("H17").Value = 46055
$ws.Range("D18").Value = 286
$ws.Range("E18").Value = 30803
$ws.Range("F18").Value = 45992
$ws.Range("H18").Value = 46055
$ws.Range("D19").Value = 424
$ws.Range("E19").Value = 30803
$ws.Range("F19").Value = 45992
$ws.Range("H19").Value = 46055
$ws.Range("E20").Value = 28460
$ws.Range("F20").Value = 45992
$ws.Range("G20").Value = 30865
$ws.Range("H20").Value = 46055
$ws.Range("C21").Value = 336
$ws.Range("F21").Value = 45992
$ws.Range("G21").Value = 30865
$ws.Range("H21").Value = 46055
$ws.Range("C22").Value = 349
$ws.Range("D22").Value = 393
$ws.Range("F22").Value = 45992
$ws.Range("H22").Value = 46055
$ws.Range("D23").Value = 341
$ws.Range("E23").Value = 30742
$ws.Range("F23").Value = 45931
$ws.Range("H23").Value = 46055
$ws.Range("D24").Value = 339
$ws.Range("H24").Value = 46055
$ws.Range("E25").Value = 30803
$ws.Range("F25").Value = 45992
$ws.Range("G25").Value = 30865
$ws.Range("H25").Value = 46055
$ws.Range("D26").Value = 406
$ws.Range("E26").Value = 30773
$ws.Range("F26").Value = 45962
$ws.Range("H26").Value = 46055
$ws.Range("D27").Value = 243
$ws.Range("H27").Value = 46055
$ws.Range("C28").Value = 420
$ws.Range("F28").Value = 45992
$ws.Range("G28").Value = 30865
$ws.Range("H28").Value = 46055
$ws.Range("E29").Value = 30803
$ws.Range("F29").Value = 45992
$ws.Range("G29").Value = 30865
$ws.Range("H29").Value = 46055
$ws.Range("C31").Value = 430
$ws.Range("D31").Value = 343
$ws.Range("F31").Value = 45962
$ws.Range("H31").Value = 46055
$ws.Range("D32").Value = 436
$ws.Range("E32").Value = 30803
$ws.Range("F32").Value = 45992
$ws.Range("H32").Value = 46055
$ws.Range("D33").Value = 343
$ws.Range("H33").Value = 46055
$ws.Range("C34").Value = 385
$ws.Range("D34").Value = 393
$ws.Range("F34").Value = 45992
$ws.Range("H34").Value = 46055
